$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 8500
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H14").Value = 8500
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H15").Value = 1282.5454
$ws.Range("I15").Value = 1282.5454
$ws.Range("K15").Value = 3847.6362
$ws.Range("M15").Value = -3678.6362
$ws.Range("H32").Value = 382.08334
$ws.Range("I32").Value = 350.25
$ws.Range("J32").Value = 398
$ws.Range("K32").Value = 350.25
$ws.Range("L32").Value = 398
$ws.Range("M32").Value = -24.25
$ws.Range("N32").Value = -1050
$ws.Range("H107").Value = 363.7619
$ws.Range("I107").Value = 384.8421
$ws.Range("J107").Value = 163.5
$ws.Range("K107").Value = 384.8421
$ws.Range("L107").Value = 163.5
$ws.Range("M107").Value = 1535.1579
$ws.Range("N107").Value = -4003.5
$ws.Range("H128").Value = 33987
$ws.Range("J128").Value = 34733.75
$ws.Range("L128").Value = 34733.75
$ws.Range("N128").Value = -44693.75
$ws.Range("H129").Value = 2291.5
$ws.Range("J129").Value = 910.73334
$ws.Range("L129").Value = 2732.20002
$ws.Range("N129").Value = -12732.20002
$ws.Range("H132").Value = 13168868
$ws.Range("I132").Value = 13168868
$ws.Range("K132").Value = 39506604
$ws.Range("M132").Value = -39504074
$ws.Range("H138").Value = 2514.6765
$ws.Range("I138").Value = 2539.1667
$ws.Range("J138").Value = 2501.318
$ws.Range("K138").Value = 7617.500100000001
$ws.Range("L138").Value = 7503.954000000001
$ws.Range("M138").Value = -2477.500100000001
$ws.Range("N138").Value = -17783.954

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1659.3
$ws.Range("I2").Value = 1784.875
$ws.Range("J2").Value = 1157
$ws.Range("K2").Value = 1784.875
$ws.Range("L2").Value = 1157
$ws.Range("M2").Value = -1671.875
$ws.Range("N2").Value = -1383
$ws.Range("H12").Value = 6000
$ws.Range("J12").Value = 6000
$ws.Range("L12").Value = 6000
$ws.Range("N12").Value = -6346
$ws.Range("H23").Value = 47628.125
$ws.Range("I23").Value = 61254.5
$ws.Range("K23").Value = 61254.5
$ws.Range("M23").Value = -60995.5
$ws.Range("H32").Value = 28754.984
$ws.Range("I32").Value = 6202.661
$ws.Range("J32").Value = 195078.38
$ws.Range("K32").Value = 6202.661
$ws.Range("L32").Value = 195078.38
$ws.Range("M32").Value = -5915.661
$ws.Range("N32").Value = -195652.38
$ws.Range("H110").Value = 125139450
$ws.Range("I110").Value = 166850930
$ws.Range("K110").Value = 166850930
$ws.Range("M110").Value = -166848885
$ws.Range("H116").Value = 1659.3
$ws.Range("I116").Value = 1784.875
$ws.Range("J116").Value = 1157
$ws.Range("K116").Value = 1784.875
$ws.Range("L116").Value = 1157
$ws.Range("M116").Value = 509.125
$ws.Range("N116").Value = -5745
$ws.Range("H122").Value = 1436.72
$ws.Range("I122").Value = 1451.3334
$ws.Range("K122").Value = 4354.0002
$ws.Range("M122").Value = -1904.0002
$ws.Range("H125").Value = 43041
$ws.Range("J125").Value = 43041
$ws.Range("L125").Value = 43041
$ws.Range("N125").Value = -52881
$ws.Range("H134").Value = 42464.5
$ws.Range("J134").Value = 42464.5
$ws.Range("L134").Value = 42464.5
$ws.Range("N134").Value = -52604.5
$ws.Range("H135").Value = 42885.43
$ws.Range("J135").Value = 42885.43
$ws.Range("L135").Value = 42885.43
$ws.Range("N135").Value = -53025.43
$ws.Range("H137").Value = 48000
$ws.Range("J137").Value = 48000
$ws.Range("L137").Value = 48000
$ws.Range("N137").Value = -58200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1659.3
$ws.Range("I3").Value = 1784.875
$ws.Range("J3").Value = 1157
$ws.Range("K3").Value = 1784.875
$ws.Range("L3").Value = 1157
$ws.Range("M3").Value = -1670.875
$ws.Range("N3").Value = -1385
$ws.Range("H22").Value = 291.5
$ws.Range("I22").Value = 309.8
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 309.8
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -136.8
$ws.Range("N22").Value = -546
$ws.Range("H86").Value = 8862.777
$ws.Range("I86").Value = 12992.546
$ws.Range("J86").Value = 2373.1428
$ws.Range("K86").Value = 12992.546
$ws.Range("L86").Value = 2373.1428
$ws.Range("M86").Value = -11869.546
$ws.Range("N86").Value = -4619.1428
$ws.Range("H89").Value = 8862.777
$ws.Range("I89").Value = 12992.546
$ws.Range("J89").Value = 2373.1428
$ws.Range("K89").Value = 64962.73
$ws.Range("L89").Value = 11865.714
$ws.Range("M89").Value = -59346.73
$ws.Range("N89").Value = -23097.714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19656.873
$ws.Range("I31").Value = 1172.2903
$ws.Range("J31").Value = 31594.834
$ws.Range("K31").Value = 1172.2903
$ws.Range("L31").Value = 31594.834
$ws.Range("M31").Value = -877.2902999999999
$ws.Range("N31").Value = -32184.834
$ws.Range("H34").Value = 19656.873
$ws.Range("I34").Value = 1172.2903
$ws.Range("J34").Value = 31594.834
$ws.Range("K34").Value = 1172.2903
$ws.Range("L34").Value = 31594.834
$ws.Range("M34").Value = -970.2902999999999
$ws.Range("N34").Value = -31998.834
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H132").Value = 136370000
$ws.Range("I132").Value = 333345730
$ws.Range("J132").Value = 62504100
$ws.Range("K132").Value = 1000037190
$ws.Range("L132").Value = 187512300
$ws.Range("M132").Value = -1000034660
$ws.Range("N132").Value = -187517360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1556.6666
$ws.Range("I26").Value = 95
$ws.Range("J26").Value = 2726
$ws.Range("K26").Value = 285
$ws.Range("L26").Value = 8178
$ws.Range("M26").Value = 3
$ws.Range("N26").Value = -8754
$ws.Range("H38").Value = 38.235294
$ws.Range("J38").Value = 56.875
$ws.Range("L38").Value = 170.625
$ws.Range("N38").Value = -864.625
$ws.Range("H58").Value = 2780
$ws.Range("I58").Value = 900
$ws.Range("K58").Value = 2700
$ws.Range("M58").Value = -2572
$ws.Range("H117").Value = 7205.222
$ws.Range("J117").Value = 7997.6875
$ws.Range("L117").Value = 23993.0625
$ws.Range("N117").Value = -30877.0625
$ws.Range("H129").Value = 13429265
$ws.Range("I129").Value = 35722908
$ws.Range("J129").Value = 424640.25
$ws.Range("K129").Value = 107168724
$ws.Range("L129").Value = 1273920.75
$ws.Range("M129").Value = -107163724
$ws.Range("N129").Value = -1283920.75
$ws.Range("H131").Value = 816.89
$ws.Range("I131").Value = 408
$ws.Range("J131").Value = 838.4105
$ws.Range("K131").Value = 1224
$ws.Range("L131").Value = 2515.2315
$ws.Range("M131").Value = 3816
$ws.Range("N131").Value = -12595.2315
$ws.Range("H132").Value = 1473.6364
$ws.Range("I132").Value = 670.5
$ws.Range("J132").Value = 9505
$ws.Range("K132").Value = 6034.5
$ws.Range("L132").Value = 85545
$ws.Range("M132").Value = -3504.5
$ws.Range("N132").Value = -90605

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 44994
$ws.Range("J64").Value = 44994
$ws.Range("L64").Value = 44994
$ws.Range("N64").Value = -45490
$ws.Range("H67").Value = 44994
$ws.Range("J67").Value = 44994
$ws.Range("L67").Value = 44994
$ws.Range("N67").Value = -46710
$ws.Range("H133").Value = 68000
$ws.Range("J133").Value = 68000
$ws.Range("L133").Value = 68000
$ws.Range("N133").Value = -78120
$ws.Range("H135").Value = 35983
$ws.Range("J135").Value = 35983
$ws.Range("L135").Value = 35983
$ws.Range("N135").Value = -46123

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 13966.286
$ws.Range("I48").Value = 11541
$ws.Range("J48").Value = 17200
$ws.Range("K48").Value = 11541
$ws.Range("L48").Value = 17200
$ws.Range("M48").Value = -10880
$ws.Range("N48").Value = -18522
$ws.Range("H132").Value = 5188.353
$ws.Range("I132").Value = 5080.2666
$ws.Range("K132").Value = 15240.7998
$ws.Range("M132").Value = -12710.7998
$ws.Range("H133").Value = 34449.125
$ws.Range("J133").Value = 34449.125
$ws.Range("L133").Value = 34449.125
$ws.Range("N133").Value = -39509.125
$ws.Range("H136").Value = 1307.659
$ws.Range("I136").Value = 1157.7567
$ws.Range("K136").Value = 3473.2701
$ws.Range("M136").Value = -923.2700999999997
$ws.Range("H137").Value = 29800
$ws.Range("J137").Value = 37000
$ws.Range("L137").Value = 37000
$ws.Range("N137").Value = -47200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 31598.857
$ws.Range("J56").Value = 36198.668
$ws.Range("L56").Value = 36198.668
$ws.Range("N56").Value = -37626.668
$ws.Range("H107").Value = 527.3333
$ws.Range("I107").Value = 390.16666
$ws.Range("J107").Value = 801.6667
$ws.Range("K107").Value = 1170.49998
$ws.Range("L107").Value = 2405.0001
$ws.Range("M107").Value = 749.5000199999999
$ws.Range("N107").Value = -6245.0001
$ws.Range("H136").Value = 2440.8235
$ws.Range("I136").Value = 677.7857
$ws.Range("J136").Value = 10668.333
$ws.Range("K136").Value = 2033.3571
$ws.Range("L136").Value = 32004.999
$ws.Range("M136").Value = 516.6428999999998
$ws.Range("N136").Value = -37104.999
